$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $found = $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $old"
    }
}

Replace-Text "107÷8=13, 3" "323÷4=80, 3"
Replace-Text "808÷5=161, 3" "747÷5=149, 2"
Replace-Text "377÷7=53, 6" "678÷6=113, 0"
Replace-Text "529÷7=75, 4" "571÷4=142, 3"
Replace-Text "298÷7=42, 4" "706÷5=141, 1"
Replace-Text "420÷4=105, 0" "238÷6=39, 4"
Replace-Text "224÷6=37, 2" "369÷3=123, 0"
Replace-Text "321÷5=64, 1" "528÷3=176, 0"
Replace-Text "888÷6=148, 0" "375÷6=62, 3"
Replace-Text "836÷9=92, 8" "442÷4=110, 2"
Replace-Text "364÷5=72, 4" "557÷2=278, 1"
Replace-Text "139÷9=15, 4" "892÷5=178, 2"
Replace-Text "342÷2=171, 0" "444÷4=111, 0"
Replace-Text "639÷4=159, 3" "204÷6=34, 0"
Replace-Text "602÷6=100, 2" "157÷9=17, 4"
Replace-Text "113÷9=12, 5" "469÷7=67, 0"
Replace-Text "577÷4=144, 1" "896÷8=112, 0"
Replace-Text "888÷9=98, 6" "783÷2=391, 1"
Replace-Text "839÷8=104, 7" "182÷3=60, 2"
Replace-Text "568÷5=113, 3" "720÷6=120, 0"
Replace-Text "779÷7=111, 2" "212÷9=23, 5"
Replace-Text "230÷5=46, 0" "837÷2=418, 1"
Replace-Text "356÷6=59, 2" "561÷3=187, 0"
Replace-Text "771÷2=385, 1" "949÷7=135, 4"
Replace-Text "537÷7=76, 5" "953÷9=105, 8"
